$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: several Price/Volume strings in this sheet look numeric to Excel's
# value-type auto-detection (e.g. "20.635.98", "1.013", "  -6.40%  "). A leading
# apostrophe forces them to be stored as literal text, exactly like a user typing
# the value manually, while leaving the cell's number format untouched.

$ws.Range("D2").Value = "'20.635.98"
$ws.Range("E2").Value = "'  -6.40%  "

$ws.Range("D3").Value = "'1.455.57"
$ws.Range("E3").Value = "'  -6.55%  "

$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "'  +1.28%  "

$ws.Range("D5").Value = "'1.010"
$ws.Range("E5").Value = "'  +1.04%  "

$ws.Range("D6").Value = "'277.68"
$ws.Range("E6").Value = "'  -4.58%  "

$ws.Range("D7").Value = "'0.3720"
$ws.Range("E7").Value = "'  -6.04%  "

$ws.Range("D8").Value = "'0.3077"
$ws.Range("E8").Value = "'  -4.64%  "

$ws.Range("D9").Value = "'41.36"
$ws.Range("E9").Value = "'  -6.59%  "

$ws.Range("B10").Value = "'Polygon"
$ws.Range("C10").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.002"
$ws.Range("E10").Value = "'  -7.21%  "

$ws.Range("B11").Value = "'Dogecoin"
$ws.Range("C11").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.06607"
$ws.Range("E11").Value = "'  -9.05%  "

$ws.Range("D12").Value = "'1.011"
$ws.Range("E12").Value = "'  +1.08%  "

$ws.Range("D13").Value = "'5.391"
$ws.Range("E13").Value = "'  -5.50%  "

$ws.Range("D14").Value = "'17.34"
$ws.Range("E14").Value = "'  -7.88%  "

$ws.Range("D15").Value = "'1.481.67"
$ws.Range("E15").Value = "'  -4.48%  "

$ws.Range("D16").Value = "'6.174"
$ws.Range("E16").Value = "'  -7.13%  "

$ws.Range("D17").Value = "'0.00001014"
$ws.Range("E17").Value = "'  -10.24%  "

$ws.Range("D18").Value = "'0.06457"
$ws.Range("E18").Value = "'  -2.02%  "

$ws.Range("D19").Value = "'78.04"
$ws.Range("E19").Value = "'  -6.71%  "

$ws.Range("D20").Value = "'1.009"
$ws.Range("E20").Value = "'  +0.93%  "

$ws.Range("D21").Value = "'5.778"
$ws.Range("E21").Value = "'  -8.03%  "

$ws.Range("D22").Value = "'14.65"
$ws.Range("E22").Value = "'  -6.01%  "

$ws.Range("D23").Value = "'10.70"
$ws.Range("E23").Value = "'  -5.71%  "

$ws.Range("D24").Value = "'2.326"
$ws.Range("E24").Value = "'  -1.85%  "

$ws.Range("D25").Value = "'20.706.89"
$ws.Range("E25").Value = "'  -6.16%  "

$ws.Range("B26").Value = "'Monero"
$ws.Range("C26").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'144.39"
$ws.Range("E26").Value = "'  -2.77%  "

$ws.Range("B27").Value = "'LidoDAOToken"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.238"
$ws.Range("E27").Value = "'  -7.88%  "

$ws.Range("D28").Value = "'17.30"
$ws.Range("E28").Value = "'  -7.19%  "

$ws.Range("D29").Value = "'4.768"
$ws.Range("E29").Value = "'  -2.16%  "

$ws.Range("D30").Value = "'1.653.36"
$ws.Range("E30").Value = "'  -4.24%  "

$ws.Range("D31").Value = "'109.66"
$ws.Range("E31").Value = "'  -8.13%  "

$ws.Range("D32").Value = "'5.547"
$ws.Range("E32").Value = "'  -5.50%  "

$ws.Range("D33").Value = "'0.9052"
$ws.Range("E33").Value = "'  -8.74%  "

$ws.Range("D34").Value = "'0.07784"
$ws.Range("E34").Value = "'  -6.48%  "

$ws.Range("D35").Value = "'8.257"
$ws.Range("E35").Value = "'  -9.93%  "

$ws.Range("D36").Value = "'1.444"
$ws.Range("E36").Value = "'  -9.87%  "

$ws.Range("B37").Value = "'Frax"
$ws.Range("C37").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.007"
$ws.Range("E37").Value = "'  +0.83%  "

$ws.Range("B38").Value = "'Aptos"
$ws.Range("C38").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'10.90"
$ws.Range("E38").Value = "'  +1.20%  "

$ws.Range("B39").Value = "'Hedera"
$ws.Range("C39").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05624"
$ws.Range("E39").Value = "'  -6.46%  "

$ws.Range("B40").Value = "'InternetComputer(DFINITY)"
$ws.Range("C40").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'4.758"
$ws.Range("E40").Value = "'  -7.41%  "

$ws.Range("D41").Value = "'0.1923"
$ws.Range("E41").Value = "'  -5.96%  "

$ws.Range("D42").Value = "'0.02046"
$ws.Range("E42").Value = "'  -9.85%  "

$ws.Range("D43").Value = "'1.113"
$ws.Range("E43").Value = "'  -7.79%  "

$ws.Range("D44").Value = "'0.5403"
$ws.Range("E44").Value = "'  -7.34%  "

$ws.Range("D45").Value = "'3.618"
$ws.Range("E45").Value = "'  -3.48%  "

$ws.Range("D46").Value = "'12.40"
$ws.Range("E46").Value = "'  -4.99%  "

$ws.Range("D47").Value = "'0.5188"
$ws.Range("E47").Value = "'  -7.29%  "

$ws.Range("D48").Value = "'1.791"
$ws.Range("E48").Value = "'  -5.97%  "

$ws.Range("D49").Value = "'109.38"
$ws.Range("E49").Value = "'  -7.66%  "

$ws.Range("D50").Value = "'1.075"
$ws.Range("E50").Value = "'  -5.70%  "

$ws.Range("D51").Value = "'0.06357"
$ws.Range("E51").Value = "'  -6.93%  "
